$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F column) values for rows 2-6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9841
$ws1.Range("F3").Value = 216
$ws1.Range("F4").Value = 40
$ws1.Range("F5").Value = 575
$ws1.Range("F6").Value = 473

# Sheet "全部类型" - update "想去人数" (F column) values for matching rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9841
$ws4.Range("F3").Value = 216
$ws4.Range("F4").Value = 40
$ws4.Range("F5").Value = 575
$ws4.Range("F7").Value = 473
